# Auto-generated edit script: updates market-price-derived columns (H-N)
# on the Leve profit sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$rows = @(
  @{ Row=19; H=1847.3; I=300; J=2019.2222; K=300; L=2019.2222; M=-125; N=-2369.2222 },
  @{ Row=28; H=73166.28999999999; I=144376.72; K=144376.72; M=-143891.72 },
  @{ Row=33; H=0; I=0; K=0; M=$null },
  @{ Row=58; H=8028.294; J=9638.357; L=28915.071; N=-29215.071 },
  @{ Row=92; H=1003.6667; I=424.83334; J=1389.5555; K=424.83334; L=1389.5555; M=823.16666; N=-3885.5555 },
  @{ Row=131; H=2916.875; I=2671; J=3457.8; K=8013; L=10373.4; M=-2973; N=-20453.4 },
  @{ Row=137; H=3343.1636; I=2076.639; K=6229.917; M=-3679.917 },
  @{ Row=138; H=4219.7095; I=2611.7407; K=7835.222099999999; M=-2695.222099999999 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("ARM")
$rows = @(
  @{ Row=32; H=3706.0847; I=3197.611; K=3197.611; M=-2910.611 },
  @{ Row=102; H=1985.32; I=1905.826; K=1905.826; M=-283.826 },
  @{ Row=111; H=62522; J=62522; L=62522; N=-70702 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("BSM")
$rows = @(
  @{ Row=99; H=3540.1765; I=3732.8333; K=3732.8333; M=-2234.8333 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("CRP")
$rows = @(
  @{ Row=3; H=4999.6665; I=0; J=4999.6665; K=0; L=4999.6665; M=$null; N=-5225.6665 },
  @{ Row=22; H=411.2; I=318.66666; K=318.66666; M=31.33334000000002 },
  @{ Row=62; H=5551.25; J=5902; L=5902; N=-7150 },
  @{ Row=65; H=5551.25; J=5902; L=29510; N=-35750 },
  @{ Row=69; H=39750; I=39750; K=39750; M=-39001 },
  @{ Row=72; H=39750; I=39750; K=119250; M=-115506 },
  @{ Row=99; H=505381.1; I=4202.75; K=4202.75; M=-2704.75 },
  @{ Row=122; H=3865.2307; I=3242.2; J=4254.625; K=9726.599999999999; L=12763.875; M=-7276.599999999999; N=-17663.875 },
  @{ Row=126; H=505381.1; I=4202.75; K=12608.25; M=-10138.25 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("CUL")
$rows = @(
  @{ Row=32; H=4714500; J=500250; L=1500750; N=-1501316 },
  @{ Row=38; H=43.8; I=52.5; K=157.5; M=189.5 },
  @{ Row=56; H=6850; I=6850; K=6850; M=-6320 },
  @{ Row=122; H=112024.78; I=1249.5; J=143674.86; K=11245.5; L=1293073.74; M=-8795.5; N=-1297973.74 },
  @{ Row=131; H=5428.8423; I=1179.4286; J=7907.6665; K=3538.2858; L=23722.9995; M=1501.7142; N=-33802.99950000001 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("GSM")
$rows = @(
  @{ Row=27; H=1000; I=1000; J=1000; K=1000; L=1000; M=-834; N=-1332 },
  @{ Row=102; H=8384.091; I=8462.105; K=8462.105; M=-6840.105 },
  @{ Row=122; H=8452; I=8091.1763; K=24273.5289; M=-21823.5289 },
  @{ Row=126; H=3221.8125; I=3036.5386; K=9109.6158; M=-6639.6158 },
  @{ Row=132; H=838134.5600000001; I=1253875.9; J=6652; K=3761627.7; L=19956; M=-3759097.7; N=-25016 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("LTW")
$rows = @(
  @{ Row=7; H=532201.2; I=914094.4399999999; K=914094.4399999999; M=-913982.4399999999 },
  @{ Row=11; H=8633.333000000001; I=11450; J=3000; K=11450; L=3000; M=-11310; N=-3280 },
  @{ Row=22; H=1250; I=1250; J=0; K=1250; L=0; M=-955; N=$null },
  @{ Row=27; H=1250; I=1250; J=0; K=1250; L=0; M=-1143; N=$null },
  @{ Row=46; H=3431.4849; I=2966.5217; K=2966.5217; M=-2778.5217 },
  @{ Row=61; H=4541.364; I=3428; J=6927.143; K=3428; L=6927.143; M=-3226; N=-7331.143 },
  @{ Row=113; H=4541.364; I=3428; J=6927.143; K=3428; L=6927.143; M=-1258; N=-11267.143 },
  @{ Row=126; H=532201.2; I=914094.4399999999; K=2742283.32; M=-2739813.32 },
  @{ Row=127; H=0; J=0; L=0; N=$null },
  @{ Row=136; H=4368.143; I=3762.8333; K=11288.4999; M=-8738.499899999999 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("WVR")
$rows = @(
  @{ Row=62; H=9330.6; I=8820; K=8820; M=-8196 },
  @{ Row=65; H=9330.6; I=8820; K=44100; M=-40980 },
  @{ Row=122; H=37041356; I=62502970; J=6281.8184; K=187508910; L=18845.4552; M=-187506460; N=-23745.4552 },
  @{ Row=126; H=3533.3076; I=3654.8572; K=10964.5716; M=-8494.571599999999 },
  @{ Row=132; H=3332.1843; I=2362.1155; J=5434; K=7086.3465; L=16302; M=-4556.3465; N=-21362 },
  @{ Row=136; H=478009.25; I=589599.6; K=1768798.8; M=-1768798.8 }
)
foreach ($r in $rows) {
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$($r.Row)"
            $val = $r[$col]
            if ($val -eq $null) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

